# Rename sample IDs so that one sample's ID is a literal substring of
# another's (test data fixture for a substring-matching bug):
#   test_sample_2_T_IGO / test_investigator_sample_2_T  -> test_sample_1a_IGO / test_investigator_sample_1a
#   test_sample_1_N_IGO / test_investigator_sample_1_N  -> test_sample_1_IGO  / test_investigator_sample_1
#
# These renames are reflected on both the "SampleInfo" sheet (rows 2 & 3)
# and the "SampleRenames" sheet (rows 2 & 3, old-name/new-name columns).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # SampleInfo
$ws2 = $wb.Worksheets.Item(2)   # SampleRenames

# --- SampleInfo ---------------------------------------------------------
$ws1.Range("C2").Value = "test_investigator_sample_1a"
$ws1.Range("C3").Value = "test_investigator_sample_1"

$ws1.Range("A2").Value = "test_sample_1a_IGO"
$ws1.Range("A3").Value = "test_sample_1_IGO"

# --- SampleRenames -------------------------------------------------------
$ws2.Range("A2").Value = "test_sample_1a_IGO"
$ws2.Range("B2").Value = "test_sample_1a"

$ws2.Range("A3").Value = "test_sample_1_IGO"
$ws2.Range("B3").Value = "test_sample_1"

# --- UI state: active tab moved from SampleInfo to SampleRenames,
#     along with each sheet's last selection -----------------------------
$ws1.Activate()
$ws1.Range("A3").Select()

$ws2.Activate()
$ws2.Range("B13").Select()
